$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.996.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.97%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.56"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.85%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.56%  "

# Row 9
$ws.Range("E9").Value = "  -3.19%  "

# Row 10
$ws.Range("E10").Value = "  -1.67%  "

# Row 11
$ws.Range("E11").Value = "  +0.33%  "

# Row 12
$ws.Range("E12").Value = "  -0.78%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.27%  "

# Row 14
$ws.Range("E14").Value = "  -1.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.639.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.948.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.120.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("E19").Value = "  +0.50%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.710"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.06%  "

# Row 23
$ws.Range("E23").Value = "  -0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "

# Row 25
$ws.Range("E25").Value = "  -3.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.50%  "

# Row 27
$ws.Range("E27").Value = "  +0.09%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "  -1.31%  "

# Row 30
$ws.Range("E30").Value = "  +0.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "

# Row 32
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
$ws.Range("E33").Value = "  -6.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.65%  "

# Row 39
$ws.Range("E39").Value = "  -1.58%  "

# Row 40
$ws.Range("E40").Value = "  -1.48%  "

# Row 41
$ws.Range("E41").Value = "  +1.15%  "

# Row 42
$ws.Range("E42").Value = "  -0.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.823.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.27%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "382.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.38%  "

# Row 45
$ws.Range("E45").Value = "  -2.06%  "

# Row 46
$ws.Range("E46").Value = "  -9.42%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.90%  "

# Row 50
$ws.Range("E50").Value = "  -0.89%  "

# Row 51
$ws.Range("E51").Value = "  -0.65%  "
